# Error Calculations and Plots
# Two records ("RM 232" and "SC 92") were removed from the missing-data
# sample entirely, shifting all subsequent rows up; the remaining rows
# also got a refreshed (re-sampled) pattern of missing ("inlineStr")
# vs. present (numeric) values in columns C/D ("B"/"C" headers) and
# D/E/F ("C"/"D"/"F" headers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two removed rows (delete the higher row index first so the
# second delete's row number is still correct).
$ws.Rows(28).Delete()
$ws.Rows(26).Delete()

# Re-apply the refreshed missing/available pattern for the remaining rows.
$ws.Range("D2").Value = -13.5
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("E5").ClearContents()
$ws.Range("F6").Value = 16.43
$ws.Range("E8").Value = -6.6
$ws.Range("E10").Value = -6.1
$ws.Range("D11").Value = -15.5
$ws.Range("F11").Value = 17.65
$ws.Range("E12").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("D13").ClearContents()
$ws.Range("F13").Value = 17.1
$ws.Range("E15").Value = -8.4
$ws.Range("F17").ClearContents()
$ws.Range("E18").ClearContents()
$ws.Range("F18").Value = 18.35
$ws.Range("E19").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("D21").Value = -14.3
$ws.Range("F24").ClearContents()
$ws.Range("D25").ClearContents()
$ws.Range("E25").Value = -7.1
$ws.Range("F25").Value = 16.6
$ws.Range("E27").Value = -10
$ws.Range("B29").ClearContents()
$ws.Range("E29").ClearContents()
$ws.Range("F31").ClearContents()
$ws.Range("F32").ClearContents()
$ws.Range("B33").Value = -19.5
$ws.Range("D33").Value = -14.1
$ws.Range("E33").ClearContents()
